try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # Insert a new column before column B, shifting existing B->C and C->D.
    $ws.Columns.Item(2).Insert()

    # Populate the new column with the person's name for each row.
    $ws.Range("B1").Value = "Демидко"
    $ws.Range("B2").Value = "ворона"

    # The inserted column picked up column A's formatting; reset it to Normal
    # so the new cells are unstyled, matching the rest of the date columns.
    $ws.Range("B1:B2").Style = "Normal"
} catch {
    Write-Output "ERROR: $_"
    throw
}
